$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36; this shifts rows 36-57 down to 37-58
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the linkeR entry
$ws.Range("A36").Value = "linkeR: Effortless Linked Views for Shiny Applications"
$ws.Range("B36").Value = "linkeR makes it effortless to create linked views in Shiny applications. When users interact with one component (like clicking a map marker), all related components (tables, charts, other maps) automatically update to show corresponding information."
$ws.Range("C36").Value = "Jake Wagoner"
$ws.Range("D36").Value = "jakew@sci.utah.edu"
$ws.Range("E36").Value = "Yes"
$ws.Range("G36").Value = "Development"
$ws.Range("H36").Value = "MIT"
$ws.Range("I36").Value = "R"
$ws.Range("J36").Value = "R Shiny Developers"
$ws.Range("K36").Value = "Moderate Programming"
$ws.Range("L36").Value = "Developer Tool"
$ws.Range("N36").Value = "https://epiforesite.github.io/linkeR/"
$ws.Range("O36").Value = "https://github.com/EpiForeSITE/linkeR/"
